$d = $word.ActiveDocument

# --- Change 1: insert a new, empty paragraph (bold formatting mark only,
# no run/text) right after the "Step 1: Connect to Azure Cloud Shell"
# heading paragraph. ---
$findRng = $d.Content
$found = $findRng.Find.Execute("Step 1: Connect to Azure Cloud Shell", $true, $false,
                                $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the 'Step 1' heading paragraph"
}
$headingPara = $findRng.Paragraphs(1)
$insertAt = $headingPara.Range
$insertAt.Collapse(0)              # wdCollapseEnd
$insertAt.InsertParagraphAfter()

$newPara = $headingPara.Next()
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
              '<w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="en-IN"/></w:rPr></w:pPr>' +
              '</w:p>'
$newPara.Range.InsertXML($newParaXml) | Out-Null

# --- Change 2: mark the run holding the page-break / "rg-powershell"
# resource-group screenshot as NoProof. That run lives in the paragraph
# immediately following the "New-AzResourceGroup ... rg-powershell ..."
# command line. ---
$cmdRng = $d.Content
$foundCmd = $cmdRng.Find.Execute('rg-powershell" -Location "South Central US"', $true,
                                  $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundCmd) {
    throw "Could not find the 'rg-powershell' command paragraph"
}
$cmdPara = $cmdRng.Paragraphs(1)
$picPara = $cmdPara.Next()
$picPara.Range.NoProofing = -1
